# The scraper now also grabs each org's revenue, so:
#   - the old "Mission" header (column J) becomes "Revenue"
#   - a new "Mission Statement" column is appended as column K
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J1").Value = "Revenue"
$ws.Range("K1").Value = "Mission Statement"

# Column J no longer holds long mission text, so its stored column
# width shrinks drastically from the old "bestFit" giant width down to
# a narrow one sized for the new Revenue numbers.
$ws.Columns.Item(10).ColumnWidth = 6.833333333333333

# Reflect where the user had scrolled/selected after adding the column.
$ws.Range("J3").Select()
